$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the entire contents of row 23 and row 24 (all columns A:AY).
# Use a temporary holding row (just past the used range) so the swap can be
# done with simple Copy operations while explicitly clearing destinations
# first (this engine's Copy does not blank out cells that are empty in the
# source, so we clear before every paste to avoid leaving stale values).

$row1 = 23
$row2 = 24
$tempRow = 26

$rSrc1 = $ws.Range("A$($row1):AY$($row1)")
$rSrc2 = $ws.Range("A$($row2):AY$($row2)")
$rTemp = $ws.Range("A$($tempRow):AY$($tempRow)")

# 1. Stash row 23 in the temp row
$rTemp.Clear()
$rSrc1.Copy($rTemp)

# 2. Move row 24 into row 23
$rSrc1.Clear()
$rSrc2.Copy($rSrc1)

# 3. Move the stashed original row 23 into row 24
$rSrc2.Clear()
$rTemp.Copy($rSrc2)

# 4. Remove the temporary row so the used range goes back to normal
$rTemp.Clear()
